$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.009.95"
$ws.Range("E2").Value = "  -3.81%  "
$ws.Range("D3").Value = "3.311.33"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.76"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.12"
$ws.Range("E6").Value = "  -3.85%  "
$ws.Range("E7").Value = "  +2.38%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.127"
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.64"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").Value = "3.889.79"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.61"
$ws.Range("E14").Value = "  -4.69%  "
$ws.Range("D15").Value = "66.141.95"
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "3.318.28"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "435.64"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.66"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.50"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.55"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.23"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.518"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000115"
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.03"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.66"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.22"
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.74"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.18"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  -4.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.62"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  -6.76%  "
$ws.Range("D39").Value = "2.830.11"
$ws.Range("E39").Value = "  +5.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.790"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.42"
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.16"
$ws.Range("E42").Value = "  -5.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.18"
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0664"
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.02"
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.33"
$ws.Range("E46").Value = "  -5.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "321.58"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.13"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.968"
$ws.Range("E51").Value = "  -3.05%  "
